$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force a hard run-boundary around [startOffset, endOffset) by
# toggling a character property on and back off again. Word (and this
# COM-interop runtime) always materialises the touched sub-range as its
# own <w:r> once the formatting differs even momentarily, and once the
# final (off) value again matches the surrounding runs' formatting it
# simply stays split into its own (identically formatted) run instead
# of re-merging. This lets us re-create a precise run layout without
# altering the visible formatting.
# ---------------------------------------------------------------------
function Isolate-Range($startOffset, $endOffset) {
    $rr = $d.Range($startOffset, $endOffset)
    $rr.Font.Bold = 1
    $rr.Font.Bold = 0
}

function Isolate-Blob($blobAbs, $relStart, $relEnd) {
    $s = $blobAbs + $relStart
    $e = $blobAbs + $relEnd
    Isolate-Range $s $e
}

# ---------------------------------------------------------------------
# Locate the paragraph that contains the sentence we need to touch.
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*De acuerdo con el*" -and $ptext -like "*artículo*") {
        $targetIndex = $i
    }
}

$para = $d.Paragraphs.Item($targetIndex).Range
$paraStart = $para.Start

# ---------------------------------------------------------------------
# Step 1: split the ", " run (right after "...Centro") into two runs,
#         "," and " " - matches the first hunk of the diff.
# ---------------------------------------------------------------------
$t0 = $para.Text
$idx0 = $t0.IndexOf(", artículo")
$spaceAbs = $paraStart + $idx0 + 1
$spaceAbsEnd = $spaceAbs + 1
Isolate-Range $spaceAbs $spaceAbsEnd

# ---------------------------------------------------------------------
# Step 2: remove the comma that sits right after "{{articulo}}".
# ---------------------------------------------------------------------
$t1 = $d.Paragraphs.Item($targetIndex).Range.Text
$idxComma1 = $t1.IndexOf("{{articulo}},") + 12
$commaAbs = $paraStart + $idxComma1
$commaAbsEnd = $commaAbs + 1
$commaRange = $d.Range($commaAbs, $commaAbsEnd)
$commaRange.Text = ""

# ---------------------------------------------------------------------
# Step 3: insert a comma right after "{{apartado}}" (i.e. before the
#         space that precedes "le comunico").
# ---------------------------------------------------------------------
$t2 = $d.Paragraphs.Item($targetIndex).Range.Text
$idxSpaceLe = $t2.IndexOf("{{apartado}} le") + 12
$spaceLeAbs = $paraStart + $idxSpaceLe
$spaceLeAbsEnd = $spaceLeAbs + 1
$spaceLeRange = $d.Range($spaceLeAbs, $spaceLeAbsEnd)
$spaceLeRange.InsertBefore(",")

# ---------------------------------------------------------------------
# Step 4: the two text edits above merged everything from "artículo "
# through to the "{{" right before the "{{nombre_alumno}}" field into
# a single run (same formatting throughout, so the engine coalesces
# touched runs). Re-carve that span back into the original run
# layout (same pieces as before, just with the comma relocated),
# matching the second hunk of the diff.
# ---------------------------------------------------------------------
$t3 = $d.Paragraphs.Item($targetIndex).Range.Text
$blobIdx = $t3.IndexOf("artículo {{articulo}}")
$blobAbs = $paraStart + $blobIdx

# relative (start,end) offsets of each piece inside the blob:
#   'artículo '              0   9
#   '{{articulo}}'           9  21
#   ' apartado {{apartado}}' 21  43
#   ','                      43  44
#   ' '                      44  45
#   'le comun'               45  53
#   'ico q'                  53  58
#   'ue su '                 58  64
#   'hij'                    64  67
#   'o/'                     67  69
#   'a'                      69  70
#   ':'                      70  71
#   ' '                      71  72
#   '{{'                     72  74  (left boundary only needed)
Isolate-Blob $blobAbs 0 9
Isolate-Blob $blobAbs 9 21
Isolate-Blob $blobAbs 21 43
Isolate-Blob $blobAbs 43 44
Isolate-Blob $blobAbs 44 45
Isolate-Blob $blobAbs 45 53
Isolate-Blob $blobAbs 53 58
Isolate-Blob $blobAbs 58 64
Isolate-Blob $blobAbs 64 67
Isolate-Blob $blobAbs 67 69
Isolate-Blob $blobAbs 69 70
Isolate-Blob $blobAbs 70 71
Isolate-Blob $blobAbs 71 72
